# Updates the cryptocurrency price/volume table on Sheet1 to the latest
# scrape (GitHub Actions data refresh). Cells hold text-formatted numbers
# (e.g. "41.473.34", "309.42") so every write forces text (@) formatting
# to stop Excel from auto-converting look-alike numbers to real numbers
# (which would also mangle trailing zeros, e.g. "89.60" -> 89.6).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($range, [string]$value) {
    # Force the cell to store $value as literal text, never a number,
    # then drop back to the workbook default style so no stray
    # number-format override is left behind on the cell.
    $range.NumberFormat = "@"
    $range.Value = $value
    $range.Style = "Normal"
}

# Row 2
$ws.Range("D2").Value = '41.473.34'
$ws.Range("E2").Value = '  -1.04%  '
# Row 3
$ws.Range("D3").Value = '2.439.84'
$ws.Range("E3").Value = '  -1.62%  '
# Row 4
Set-TextValue $ws.Range("D4") '1.01'
$ws.Range("E4").Value = '  +1.14%  '
# Row 5
Set-TextValue $ws.Range("D5") '309.42'
$ws.Range("E5").Value = '  -0.77%  '
# Row 6
Set-TextValue $ws.Range("D6") '89.60'
$ws.Range("E6").Value = '  -5.98%  '
# Row 7
Set-TextValue $ws.Range("D7") '0.530'
$ws.Range("E7").Value = '  -4.39%  '
# Row 8
$ws.Range("E8").Value = '  +0.94%  '
# Row 9
Set-TextValue $ws.Range("D9") '0.481'
$ws.Range("E9").Value = '  -6.16%  '
# Row 10
Set-TextValue $ws.Range("D10") '31.61'
$ws.Range("E10").Value = '  -7.23%  '
# Row 11
Set-TextValue $ws.Range("D11") '0.0765'
$ws.Range("E11").Value = '  -2.80%  '
# Row 12
Set-TextValue $ws.Range("D12") '0.109'
$ws.Range("E12").Value = '  +0.01%  '
# Row 13
$ws.Range("D13").Value = '2.811.47'
$ws.Range("E13").Value = '  -1.61%  '
# Row 14
Set-TextValue $ws.Range("D14") '6.67'
$ws.Range("E14").Value = '  -5.24%  '
# Row 15
$ws.Range("B15").Value = 'WrappedEther'
$ws.Range("C15").Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range("D15").Value = '2.496.01'
$ws.Range("E15").Value = '  +0.74%  '
# Row 16
$ws.Range("B16").Value = 'Chainlink'
$ws.Range("C16").Value = 'https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link'
Set-TextValue $ws.Range("D16") '14.93'
$ws.Range("E16").Value = '  +0.65%  '
# Row 17
Set-TextValue $ws.Range("D17") '0.749'
$ws.Range("E17").Value = '  -5.38%  '
# Row 18
$ws.Range("D18").Value = '41.113.22'
$ws.Range("E18").Value = '  -1.80%  '
# Row 19
Set-TextValue $ws.Range("D19") '6.13'
$ws.Range("E19").Value = '  -4.46%  '
# Row 20
$ws.Range("D20").Value = '0.0₃0898'
$ws.Range("E20").Value = '  -2.62%  '
# Row 21
Set-TextValue $ws.Range("D21") '68.76'
$ws.Range("E21").Value = '  -0.43%  '
# Row 22
Set-TextValue $ws.Range("D22") '10.67'
$ws.Range("E22").Value = '  -9.48%  '
# Row 23
Set-TextValue $ws.Range("D23") '230.43'
$ws.Range("E23").Value = '  -3.03%  '
# Row 24
Set-TextValue $ws.Range("D24") '2.66'
$ws.Range("E24").Value = '  -5.13%  '
# Row 25
$ws.Range("E25").Value = '  +0.11%  '
# Row 26
Set-TextValue $ws.Range("D26") '1.84'
$ws.Range("E26").Value = '  -5.51%  '
# Row 27
Set-TextValue $ws.Range("D27") '23.45'
$ws.Range("E27").Value = '  -5.37%  '
# Row 28
Set-TextValue $ws.Range("D28") '2.21'
$ws.Range("E28").Value = '  -0.66%  '
# Row 29
Set-TextValue $ws.Range("D29") '9.43'
$ws.Range("E29").Value = '  -3.49%  '
# Row 30
Set-TextValue $ws.Range("D30") '34.73'
$ws.Range("E30").Value = '  -5.20%  '
# Row 31
Set-TextValue $ws.Range("D31") '151.14'
$ws.Range("E31").Value = '  -2.33%  '
# Row 32
Set-TextValue $ws.Range("D32") '5.22'
$ws.Range("E32").Value = '  -7.83%  '
# Row 33
$ws.Range("E33").Value = '  -3.76%  '
# Row 34
$ws.Range("B34").Value = 'Hedera'
$ws.Range("C34").Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
Set-TextValue $ws.Range("D34") '0.0741'
$ws.Range("E34").Value = '  -2.53%  '
# Row 35
$ws.Range("B35").Value = 'ApeXProtocol'
$ws.Range("C35").Value = 'https://coinranking.com/coin/ze0N2Rcyu+apexprotocol-apex'
Set-TextValue $ws.Range("D35") '2.49'
$ws.Range("E35").Value = '  -0.40%  '
# Row 36
Set-TextValue $ws.Range("D36") '17.35'
$ws.Range("E36").Value = '  +0.66%  '
# Row 37
Set-TextValue $ws.Range("D37") '2.88'
$ws.Range("E37").Value = '  -5.12%  '
# Row 38
Set-TextValue $ws.Range("D38") '1.77'
$ws.Range("E38").Value = '  -6.46%  '
# Row 39
Set-TextValue $ws.Range("D39") '0.111'
$ws.Range("E39").Value = '  -3.68%  '
# Row 40
Set-TextValue $ws.Range("D40") '0.0980'
$ws.Range("E40").Value = '  -8.72%  '
# Row 41
Set-TextValue $ws.Range("D41") '3.98'
$ws.Range("E41").Value = '  -1.72%  '
# Row 42
$ws.Range("E42").Value = '  +1.40%  '
# Row 43
Set-TextValue $ws.Range("D43") '18.95'
$ws.Range("E43").Value = '  -11.23%  '
# Row 44
$ws.Range("D44").Value = '1.912.28'
$ws.Range("E44").Value = '  -4.61%  '
# Row 45
Set-TextValue $ws.Range("D45") '0.0273'
$ws.Range("E45").Value = '  -4.97%  '
# Row 46
Set-TextValue $ws.Range("D46") '2.84'
$ws.Range("E46").Value = '  -8.33%  '
# Row 47
$ws.Range("B47").Value = 'FraxShare'
$ws.Range("C47").Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
Set-TextValue $ws.Range("D47") '8.54'
$ws.Range("E47").Value = '  -1.80%  '
# Row 48
$ws.Range("B48").Value = 'RocketPoolETH'
$ws.Range("C48").Value = 'https://coinranking.com/coin/QJZRUGyNI+rocketpooleth-reth'
$ws.Range("D48").Value = '2.671.31'
$ws.Range("E48").Value = '  -1.68%  '
# Row 49
Set-TextValue $ws.Range("D49") '93.07'
$ws.Range("E49").Value = '  -5.26%  '
# Row 50
Set-TextValue $ws.Range("D50") '0.171'
$ws.Range("E50").Value = '  -6.42%  '
# Row 51
Set-TextValue $ws.Range("D51") '64.80'
$ws.Range("E51").Value = '  -7.52%  '
